$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8729138970375061
$ws.Range("B1").Value = 1.554526805877686
$ws.Range("C1").Value = 6.376989364624023
$ws.Range("D1").Value = 3.006787300109863
$ws.Range("E1").Value = 1.563802242279053
